{"js": "// Update the two-digit-divided-by-one-digit practice table: replace each\n// equation's text with a newly generated problem, cell by cell, in document\n// order. Formatting (font/size/alignment) is preserved because we replace\n// only the text inside the existing paragraph range rather than rebuilding\n// the cell body.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Map of [rowIndex, colIndex] -> new text, rowIndex/colIndex are 0-based\n// positions in the Word table grid (blank spacer rows included), matching\n// the document's row layout: data lives in rows 0, 4, 8, 12, 16.\nconst replacements = [\n  [0, 0, \"62\u00f74=\"],\n  [0, 1, \"38\u00f76=\"],\n  [0, 2, \"84\u00f75=\"],\n  [0, 3, \"55\u00f77=\"],\n  [0, 4, \"43\u00f74=\"],\n\n  [4, 0, \"66\u00f77=\"],\n  [4, 1, \"31\u00f73=\"],\n  [4, 2, \"54\u00f73=\"],\n  [4, 3, \"99\u00f72=\"],\n  [4, 4, \"48\u00f79=\"],\n\n  [8, 0, \"38\u00f77=\"],\n  [8, 1, \"16\u00f76=\"],\n  [8, 2, \"61\u00f77=\"],\n  [8, 3, \"81\u00f78=\"],\n  [8, 4, \"21\u00f75=\"],\n\n  [12, 0, \"33\u00f78=\"],\n  [12, 1, \"36\u00f79=\"],\n  [12, 2, \"88\u00f77=\"],\n  [12, 3, \"25\u00f77=\"],\n  [12, 4, \"59\u00f76=\"],\n\n  [16, 0, \"39\u00f77=\"],\n  [16, 1, \"90\u00f75=\"],\n  [16, 2, \"11\u00f79=\"],\n  [16, 3, \"12\u00f75=\"],\n  [16, 4, \"61\u00f72=\"],\n];\n\nfor (const [row, col, text] of replacements) {\n  const cell = table.getCell(row, col);\n  const range = cell.body.getRange();\n  range.insertText(text, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the two-digit-divided-by-one-digit practice table: replace each\n# equation's text with a newly generated problem, cell by cell. Writing\n# through Cell.Range.Text only swaps the run's text (leaving paragraph /\n# run formatting such as font, size and alignment untouched).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Word COM collections are 1-based. The data rows in the table are rows\n# 1, 5, 9, 13, 17 (the rows in between are blank spacer rows), each with\n# 5 columns.\n$replacements = @(\n    @{Row=1;  Col=1; Text=\"62\u00f74=\"},\n    @{Row=1;  Col=2; Text=\"38\u00f76=\"},\n    @{Row=1;  Col=3; Text=\"84\u00f75=\"},\n    @{Row=1;  Col=4; Text=\"55\u00f77=\"},\n    @{Row=1;  Col=5; Text=\"43\u00f74=\"},\n\n    @{Row=5;  Col=1; Text=\"66\u00f77=\"},\n    @{Row=5;  Col=2; Text=\"31\u00f73=\"},\n    @{Row=5;  Col=3; Text=\"54\u00f73=\"},\n    @{Row=5;  Col=4; Text=\"99\u00f72=\"},\n    @{Row=5;  Col=5; Text=\"48\u00f79=\"},\n\n    @{Row=9;  Col=1; Text=\"38\u00f77=\"},\n    @{Row=9;  Col=2; Text=\"16\u00f76=\"},\n    @{Row=9;  Col=3; Text=\"61\u00f77=\"},\n    @{Row=9;  Col=4; Text=\"81\u00f78=\"},\n    @{Row=9;  Col=5; Text=\"21\u00f75=\"},\n\n    @{Row=13; Col=1; Text=\"33\u00f78=\"},\n    @{Row=13; Col=2; Text=\"36\u00f79=\"},\n    @{Row=13; Col=3; Text=\"88\u00f77=\"},\n    @{Row=13; Col=4; Text=\"25\u00f77=\"},\n    @{Row=13; Col=5; Text=\"59\u00f76=\"},\n\n    @{Row=17; Col=1; Text=\"39\u00f77=\"},\n    @{Row=17; Col=2; Text=\"90\u00f75=\"},\n    @{Row=17; Col=3; Text=\"11\u00f79=\"},\n    @{Row=17; Col=4; Text=\"12\u00f75=\"},\n    @{Row=17; Col=5; Text=\"61\u00f72=\"}\n)\n\nforeach ($r in $replacements) {\n    $cell = $t.Cell($r.Row, $r.Col)\n    $cell.Range.Text = $r.Text\n}\n"}
